$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# isa picked up / is now working on validação_Resposta() -> mark the Check cell (merged C16:C18)
# as "in progress" (yellow, bold) and write her name in it.
$ws.Range("C16").Value = "isa"
$ws.Range("C16:C18").Interior.Color = 65535
$ws.Range("C16:C18").Font.Bold = $true

# Clear the stray "Em análise" note that used to live in D20 and reset its formatting
# back to a plain (white, borderless) cell.
$ws.Range("D20").ClearContents()
$ws.Range("D20").Interior.ThemeColor = 2
$ws.Range("D20").Borders.LineStyle = -4142

# D16 becomes an (empty) placeholder cell styled like the other "note" cells (e.g. F12).
$ws.Range("D16").Font.Underline = 2

# Select D16 to match the saved cursor position.
$ws.Range("D16").Select()
